$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 totals (consolidated export order)
$ws.Range("I3").Value = 12
$ws.Range("K3").Value = 200000
$ws.Range("N3").Value = 2936000

# Remove the now-superfluous rows 4 and 5 (their data was folded into row 3)
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
